$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# Insert a new row above the current row 6 (the empty, styled row),
# shifting it down to row 7, and fill the new row 6 with the
# "style" / "default" key-value pair, matching the other rows'
# formatting (bold/orange key in column A).
$ws.Rows("6:6").Insert()

$ws.Cells.Item(6, 1).Value = "style"
$ws.Cells.Item(6, 2).Value = "default"
